# Update cryptocurrency price (D) and volume-change (E) columns
# to the latest scraped values. D-column cells that read as plain
# numeric literals need NumberFormat forced to text first so Excel
# stores them as the exact literal string (matching the original
# inlineStr cells) instead of re-parsing/rounding them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.905.78"
$ws.Range("E2").Value = "  +2.69%  "

$ws.Range("D3").Value = "2.951.99"
$ws.Range("E3").Value = "  +0.72%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.84"
$ws.Range("E5").Value = "  -0.62%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.09"
$ws.Range("E6").Value = "  +2.82%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("D8").Value = "2.948.23"
$ws.Range("E8").Value = "  +0.72%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.504"
$ws.Range("E9").Value = "  +0.68%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.10"
$ws.Range("E10").Value = "  +1.85%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.149"
$ws.Range("E11").Value = "  +5.71%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.439"
$ws.Range("E12").Value = "  +0.60%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000233"
$ws.Range("E13").Value = "  +4.57%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.68"
$ws.Range("E14").Value = "  -1.67%  "

$ws.Range("E15").Value = "  -0.66%  "

$ws.Range("D16").Value = "3.440.27"
$ws.Range("E16").Value = "  +0.82%  "

$ws.Range("D17").Value = "62.894.97"
$ws.Range("E17").Value = "  +2.75%  "

$ws.Range("E18").Value = "  +0.23%  "

$ws.Range("D19").Value = "2.953.05"
$ws.Range("E19").Value = "  +0.93%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "440.21"
$ws.Range("E20").Value = "  +1.67%  "

$ws.Range("E21").Value = "  -0.24%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.667"
$ws.Range("E22").Value = "  -0.76%  "

$ws.Range("E23").Value = "  -0.56%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.22"
$ws.Range("E24").Value = "  +2.97%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.81"
$ws.Range("E25").Value = "  -1.00%  "

$ws.Range("E26").Value = "  -1.72%  "

$ws.Range("E27").Value = "  +0.71%  "

$ws.Range("E28").Value = "  +0.02%  "

$ws.Range("E29").Value = "  +1.48%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.30"
$ws.Range("E30").Value = "  +5.79%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.60"
$ws.Range("E31").Value = "  +0.53%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0000102"
$ws.Range("E32").Value = "  +17.14%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.36"
$ws.Range("E33").Value = "  -0.73%  "

$ws.Range("E34").Value = "  -1.24%  "

$ws.Range("E35").Value = "  +0.16%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.990"
$ws.Range("E36").Value = "  -2.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.61"
$ws.Range("E37").Value = "  -0.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.06"
$ws.Range("E38").Value = "  +3.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.72"
$ws.Range("E39").Value = "  -0.10%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.03"
$ws.Range("E40").Value = "  +2.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.48"
$ws.Range("E41").Value = "  -0.48%  "

$ws.Range("E42").Value = "  -3.37%  "

$ws.Range("E43").Value = "  -0.08%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.60"
$ws.Range("E44").Value = "  -6.09%  "

$ws.Range("D45").Value = "2.708.09"
$ws.Range("E45").Value = "  +0.63%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "135.72"
$ws.Range("E46").Value = "  +1.62%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0338"
$ws.Range("E47").Value = "  -1.60%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "360.90"
$ws.Range("E48").Value = "  -0.66%  "

$ws.Range("E50").Value = "  -0.40%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.74"
$ws.Range("E51").Value = "  -3.41%  "

